$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 141-156 (Bathroom, No Motion, Inactive) ---
$wsPIR = $wb.Worksheets.Item("PIR")
$wsPIR.Range("A141").Value = "'2026-01-30"
$wsPIR.Range("B141").Value = "17:07:57"
$wsPIR.Range("C141").Value = "17:00"
$wsPIR.Range("D141").Value = "Bathroom"
$wsPIR.Range("E141").Value = "No Motion"
$wsPIR.Range("F141").Value = "Inactive"

$wsPIR.Range("A142").Value = "'2026-01-30"
$wsPIR.Range("B142").Value = "17:07:57"
$wsPIR.Range("C142").Value = "17:00"
$wsPIR.Range("D142").Value = "Bathroom"
$wsPIR.Range("E142").Value = "No Motion"
$wsPIR.Range("F142").Value = "Inactive"

$wsPIR.Range("A143").Value = "'2026-01-30"
$wsPIR.Range("B143").Value = "17:07:58"
$wsPIR.Range("C143").Value = "17:00"
$wsPIR.Range("D143").Value = "Bathroom"
$wsPIR.Range("E143").Value = "No Motion"
$wsPIR.Range("F143").Value = "Inactive"

$wsPIR.Range("A144").Value = "'2026-01-30"
$wsPIR.Range("B144").Value = "17:07:59"
$wsPIR.Range("C144").Value = "17:00"
$wsPIR.Range("D144").Value = "Bathroom"
$wsPIR.Range("E144").Value = "No Motion"
$wsPIR.Range("F144").Value = "Inactive"

$wsPIR.Range("A145").Value = "'2026-01-30"
$wsPIR.Range("B145").Value = "17:08:00"
$wsPIR.Range("C145").Value = "17:00"
$wsPIR.Range("D145").Value = "Bathroom"
$wsPIR.Range("E145").Value = "No Motion"
$wsPIR.Range("F145").Value = "Inactive"

$wsPIR.Range("A146").Value = "'2026-01-30"
$wsPIR.Range("B146").Value = "17:08:02"
$wsPIR.Range("C146").Value = "17:00"
$wsPIR.Range("D146").Value = "Bathroom"
$wsPIR.Range("E146").Value = "No Motion"
$wsPIR.Range("F146").Value = "Inactive"

$wsPIR.Range("A147").Value = "'2026-01-30"
$wsPIR.Range("B147").Value = "17:08:07"
$wsPIR.Range("C147").Value = "17:00"
$wsPIR.Range("D147").Value = "Bathroom"
$wsPIR.Range("E147").Value = "No Motion"
$wsPIR.Range("F147").Value = "Inactive"

$wsPIR.Range("A148").Value = "'2026-01-30"
$wsPIR.Range("B148").Value = "17:08:12"
$wsPIR.Range("C148").Value = "17:00"
$wsPIR.Range("D148").Value = "Bathroom"
$wsPIR.Range("E148").Value = "No Motion"
$wsPIR.Range("F148").Value = "Inactive"

$wsPIR.Range("A149").Value = "'2026-01-30"
$wsPIR.Range("B149").Value = "17:08:17"
$wsPIR.Range("C149").Value = "17:00"
$wsPIR.Range("D149").Value = "Bathroom"
$wsPIR.Range("E149").Value = "No Motion"
$wsPIR.Range("F149").Value = "Inactive"

$wsPIR.Range("A150").Value = "'2026-01-30"
$wsPIR.Range("B150").Value = "17:08:22"
$wsPIR.Range("C150").Value = "17:00"
$wsPIR.Range("D150").Value = "Bathroom"
$wsPIR.Range("E150").Value = "No Motion"
$wsPIR.Range("F150").Value = "Inactive"

$wsPIR.Range("A151").Value = "'2026-01-30"
$wsPIR.Range("B151").Value = "17:08:27"
$wsPIR.Range("C151").Value = "17:00"
$wsPIR.Range("D151").Value = "Bathroom"
$wsPIR.Range("E151").Value = "No Motion"
$wsPIR.Range("F151").Value = "Inactive"

$wsPIR.Range("A152").Value = "'2026-01-30"
$wsPIR.Range("B152").Value = "17:08:32"
$wsPIR.Range("C152").Value = "17:00"
$wsPIR.Range("D152").Value = "Bathroom"
$wsPIR.Range("E152").Value = "No Motion"
$wsPIR.Range("F152").Value = "Inactive"

$wsPIR.Range("A153").Value = "'2026-01-30"
$wsPIR.Range("B153").Value = "17:08:37"
$wsPIR.Range("C153").Value = "17:00"
$wsPIR.Range("D153").Value = "Bathroom"
$wsPIR.Range("E153").Value = "No Motion"
$wsPIR.Range("F153").Value = "Inactive"

$wsPIR.Range("A154").Value = "'2026-01-30"
$wsPIR.Range("B154").Value = "17:08:42"
$wsPIR.Range("C154").Value = "17:00"
$wsPIR.Range("D154").Value = "Bathroom"
$wsPIR.Range("E154").Value = "No Motion"
$wsPIR.Range("F154").Value = "Inactive"

$wsPIR.Range("A155").Value = "'2026-01-30"
$wsPIR.Range("B155").Value = "17:08:47"
$wsPIR.Range("C155").Value = "17:00"
$wsPIR.Range("D155").Value = "Bathroom"
$wsPIR.Range("E155").Value = "No Motion"
$wsPIR.Range("F155").Value = "Inactive"

$wsPIR.Range("A156").Value = "'2026-01-30"
$wsPIR.Range("B156").Value = "17:08:52"
$wsPIR.Range("C156").Value = "17:00"
$wsPIR.Range("D156").Value = "Bathroom"
$wsPIR.Range("E156").Value = "No Motion"
$wsPIR.Range("F156").Value = "Inactive"

# --- Humidity sheet: append rows 88-101 (Bathroom, percentage, Active) ---
$wsHumidity = $wb.Worksheets.Item("Humidity")
$wsHumidity.Range("A88").Value = "'2026-01-30"
$wsHumidity.Range("B88").Value = "17:07:57"
$wsHumidity.Range("C88").Value = "17:00"
$wsHumidity.Range("D88").Value = "Bathroom"
$wsHumidity.Range("E88").Value = "'87.4%"
$wsHumidity.Range("F88").Value = "Active"

$wsHumidity.Range("A89").Value = "'2026-01-30"
$wsHumidity.Range("B89").Value = "17:07:58"
$wsHumidity.Range("C89").Value = "17:00"
$wsHumidity.Range("D89").Value = "Bathroom"
$wsHumidity.Range("E89").Value = "'87.5%"
$wsHumidity.Range("F89").Value = "Active"

$wsHumidity.Range("A90").Value = "'2026-01-30"
$wsHumidity.Range("B90").Value = "17:07:58"
$wsHumidity.Range("C90").Value = "17:00"
$wsHumidity.Range("D90").Value = "Bathroom"
$wsHumidity.Range("E90").Value = "'87.5%"
$wsHumidity.Range("F90").Value = "Active"

$wsHumidity.Range("A91").Value = "'2026-01-30"
$wsHumidity.Range("B91").Value = "17:07:59"
$wsHumidity.Range("C91").Value = "17:00"
$wsHumidity.Range("D91").Value = "Bathroom"
$wsHumidity.Range("E91").Value = "'87.5%"
$wsHumidity.Range("F91").Value = "Active"

$wsHumidity.Range("A92").Value = "'2026-01-30"
$wsHumidity.Range("B92").Value = "17:08:02"
$wsHumidity.Range("C92").Value = "17:00"
$wsHumidity.Range("D92").Value = "Bathroom"
$wsHumidity.Range("E92").Value = "'87.4%"
$wsHumidity.Range("F92").Value = "Active"

$wsHumidity.Range("A93").Value = "'2026-01-30"
$wsHumidity.Range("B93").Value = "17:08:07"
$wsHumidity.Range("C93").Value = "17:00"
$wsHumidity.Range("D93").Value = "Bathroom"
$wsHumidity.Range("E93").Value = "'87.4%"
$wsHumidity.Range("F93").Value = "Active"

$wsHumidity.Range("A94").Value = "'2026-01-30"
$wsHumidity.Range("B94").Value = "17:08:12"
$wsHumidity.Range("C94").Value = "17:00"
$wsHumidity.Range("D94").Value = "Bathroom"
$wsHumidity.Range("E94").Value = "'87.4%"
$wsHumidity.Range("F94").Value = "Active"

$wsHumidity.Range("A95").Value = "'2026-01-30"
$wsHumidity.Range("B95").Value = "17:08:17"
$wsHumidity.Range("C95").Value = "17:00"
$wsHumidity.Range("D95").Value = "Bathroom"
$wsHumidity.Range("E95").Value = "'87.4%"
$wsHumidity.Range("F95").Value = "Active"

$wsHumidity.Range("A96").Value = "'2026-01-30"
$wsHumidity.Range("B96").Value = "17:08:22"
$wsHumidity.Range("C96").Value = "17:00"
$wsHumidity.Range("D96").Value = "Bathroom"
$wsHumidity.Range("E96").Value = "'87.4%"
$wsHumidity.Range("F96").Value = "Active"

$wsHumidity.Range("A97").Value = "'2026-01-30"
$wsHumidity.Range("B97").Value = "17:08:27"
$wsHumidity.Range("C97").Value = "17:00"
$wsHumidity.Range("D97").Value = "Bathroom"
$wsHumidity.Range("E97").Value = "'87.4%"
$wsHumidity.Range("F97").Value = "Active"

$wsHumidity.Range("A98").Value = "'2026-01-30"
$wsHumidity.Range("B98").Value = "17:08:32"
$wsHumidity.Range("C98").Value = "17:00"
$wsHumidity.Range("D98").Value = "Bathroom"
$wsHumidity.Range("E98").Value = "'87.4%"
$wsHumidity.Range("F98").Value = "Active"

$wsHumidity.Range("A99").Value = "'2026-01-30"
$wsHumidity.Range("B99").Value = "17:08:37"
$wsHumidity.Range("C99").Value = "17:00"
$wsHumidity.Range("D99").Value = "Bathroom"
$wsHumidity.Range("E99").Value = "'86.5%"
$wsHumidity.Range("F99").Value = "Active"

$wsHumidity.Range("A100").Value = "'2026-01-30"
$wsHumidity.Range("B100").Value = "17:08:42"
$wsHumidity.Range("C100").Value = "17:00"
$wsHumidity.Range("D100").Value = "Bathroom"
$wsHumidity.Range("E100").Value = "'87.4%"
$wsHumidity.Range("F100").Value = "Active"

$wsHumidity.Range("A101").Value = "'2026-01-30"
$wsHumidity.Range("B101").Value = "17:08:52"
$wsHumidity.Range("C101").Value = "17:00"
$wsHumidity.Range("D101").Value = "Bathroom"
$wsHumidity.Range("E101").Value = "'87.4%"
$wsHumidity.Range("F101").Value = "Active"

# --- mmWave sheet: append rows 39-40 (Living Room, fall detection emergency) ---
$wsMmWave = $wb.Worksheets.Item("mmWave")
$wsMmWave.Range("A39").Value = "'2026-01-30"
$wsMmWave.Range("B39").Value = "17:08:00"
$wsMmWave.Range("C39").Value = "17:00"
$wsMmWave.Range("D39").Value = "Living Room"
$wsMmWave.Range("E39").Value = "FALL_DETECTED"
$wsMmWave.Range("F39").Value = "EMERGENCY"

$wsMmWave.Range("A40").Value = "'2026-01-30"
$wsMmWave.Range("B40").Value = "17:08:40"
$wsMmWave.Range("C40").Value = "17:00"
$wsMmWave.Range("D40").Value = "Living Room"
$wsMmWave.Range("E40").Value = "PRESENCE_DETECTED"
$wsMmWave.Range("F40").Value = "Active"

Write-Output "Applied sensor log updates to PIR, Humidity, and mmWave sheets."
